$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "BUSY" legend label at M15:O15 (mirrors the PAY FLIPFLOP legend box at L16:P16) ---
$ws.Range("L16").Copy()
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("M16").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("M15").Value = "BUSY"

# --- Header row 9: swap I/C headers, drop the "P" column, add BUSY (V/R/F/B) columns ---
$ws.Range("H9").Value = "I"
$ws.Range("I9").Value = "C"
$ws.Range("J9").Value = "V"
$ws.Range("K9").Value = "R"
$ws.Range("L9").Value = "F"
$ws.Range("M9").Value = "B"

# --- New "BUSY" block mirrors the existing M10 data cell + the M11:M14 thin-left borders ---
$ws.Range("L10").Copy()
$ws.Range("M10").PasteSpecial(-4122)

$ws.Range("L11").Copy()
$ws.Range("M11").PasteSpecial(-4122)

$ws.Range("L12").Copy()
$ws.Range("M12").PasteSpecial(-4122)

$ws.Range("L13").Copy()
$ws.Range("M13").PasteSpecial(-4122)

$ws.Range("L14").Copy()
$ws.Range("M14").PasteSpecial(-4122)

# --- Selection moves to M9 ---
$ws.Range("M9").Select()

# --- Minimize the workbook window ---
$excel.ActiveWindow.WindowState = -4140

$excel.CutCopyMode = $false

"ok"
